# Generate Report for Handback
# Refresh the timestamps recorded on the "Overview" summary sheet and on the
# per-locale ("zh-cn" / "de-de") detail sheets to reflect the latest
# handoff/handback run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-29 03:05:03"

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-29 03:04:57"
$zhcn.Range("K2").Value = "2016-08-29 03:05:43"

# --- de-de detail sheet -----------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-29 03:05:03"
$dede.Range("K2").Value = "2016-08-29 03:05:49"
